$p = $ppt.ActivePresentation
$n = $p.Slides.InsertFromFile("/tmp/work/before.pptx", 1, 0, -1)
Write-Output ("inserted: " + $n)
